# Rename the "Label" column header to "What" on Sheet1 (cell B1).
# This reflects the commit's "added location parameter" change: the
# Fill scenario's second column, previously labelled "Label", now
# documents the element-location parameter and is headed "What".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "What"

# Leave the active selection on the edited cell, matching the
# post-edit cursor position.
[void]$ws.Range("B1").Select()
